$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data row (row 2) with new match values
$ws.Range("A2").Value = "KKku3143"
$ws.Range("C2").Value = "09:00"
$ws.Range("D2").Value = "INDONESIA - LIGA 1"
$ws.Range("E2").Value = "Persija Jakarta"
$ws.Range("F2").Value = "Madura United"
$ws.Range("G2").Value = 1.55
$ws.Range("H2").Value = 3.8
$ws.Range("I2").Value = 5.5
$ws.Range("J2").Value = 2.1
$ws.Range("K2").Value = 2.2
$ws.Range("L2").Value = 5.4
$ws.Range("M2").Value = 1.01
$ws.Range("N2").Value = 7.8
$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 3.2
$ws.Range("Q2").Value = 1.8
$ws.Range("R2").Value = 1.9
$ws.Range("S2").Value = 1.39
$ws.Range("T2").Value = 2.57
$ws.Range("U2").Value = 1.82
$ws.Range("W2").Value = 6.6
$ws.Range("X2").Value = 7.1
$ws.Range("Y2").Value = 8
$ws.Range("Z2").Value = 11
$ws.Range("AA2").Value = 12.5
$ws.Range("AC2").Value = 11
$ws.Range("AD2").Value = 7.5
$ws.Range("AE2").Value = 16.5
$ws.Range("AF2").Value = 75
$ws.Range("AG2").Value = 600
$ws.Range("AH2").Value = 15.5
$ws.Range("AI2").Value = 35
$ws.Range("AJ2").Value = 17.5
$ws.Range("AK2").Value = 110
$ws.Range("AL2").Value = 55
$ws.Range("AM2").Value = 55
$ws.Range("AN2").Value = 3.35
$ws.Range("AO2").Value = 7.4
$ws.Range("AQ2").Value = 23
$ws.Range("AR2").Value = 55
$ws.Range("AS2").Value = 250
$ws.Range("AT2").Value = 2.62
$ws.Range("AU2").Value = 7.6
$ws.Range("AV2").Value = 70
$ws.Range("AW2").Value = 6.9
$ws.Range("AX2").Value = 30
$ws.Range("AY2").Value = 35
$ws.Range("AZ2").Value = 200
$ws.Range("BA2").Value = 200
$ws.Range("BB2").Value = 450

# Remove the trailing Odd_CS_3-3_HT / Odd_CS_4-4_HT columns (BC:BD)
$ws.Range("BC1:BD2").EntireColumn.Delete()
